$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a flat weekly price-history table (columns A:R, header row 1).
# The edit inserts two new weekly observation rows:
#   - one new row at (original) position 244
#   - one new row further down, which ends up at position 282 once the
#     first insertion has shifted everything below it down by one
# Both new rows are duplicates of the very first data row in this block
# (old row 244: K=3500, L=3500, M=3500, N="$/docena de atados",
#  O="Region Metropolitana", P=1167, plus the shared A/B/C/E/F/G/H/I/Q/R
#  boilerplate for this market/category), with only the date (D) and,
# for the first one, the volume (J) edited to new values.
# ---------------------------------------------------------------------------

function Set-PriceRow {
    param($row, $date, $volumen)

    $ws.Cells.Item($row, 1).Value  = 4
    $ws.Cells.Item($row, 2).Value  = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($row, 3).Value  = "Los Lagos"
    $ws.Cells.Item($row, 4).Value  = $date
    $ws.Cells.Item($row, 5).Value  = 10
    $ws.Cells.Item($row, 6).Value  = 100112039
    $ws.Cells.Item($row, 7).Value  = "Ciboulette"
    $ws.Cells.Item($row, 8).Value  = "Sin especificar"
    $ws.Cells.Item($row, 9).Value  = "Primera"
    $ws.Cells.Item($row, 10).Value = $volumen
    $ws.Cells.Item($row, 11).Value = 3500
    $ws.Cells.Item($row, 12).Value = 3500
    $ws.Cells.Item($row, 13).Value = 3500
    $ws.Cells.Item($row, 14).Value = "`$/docena de atados"
    $ws.Cells.Item($row, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($row, 16).Value = 1167
    $ws.Cells.Item($row, 17).Value = 3
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}

# Insert the first new row at position 244 (everything from old row 244
# downward shifts down by one row).
$ws.Rows(244).Insert()
Set-PriceRow 244 45006 240

# Insert the second new row. Before this insert, the row that used to be
# old row 281 now sits at position 282 (because of the first insertion
# above); inserting here pushes it (and everything below it) down by one
# more row, landing the new row at position 282.
$ws.Rows(282).Insert()
Set-PriceRow 282 45005 80
